# Auto-generated edit script: updates Leve profit-calculation cells
# across the 8 crafting-job sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR)
# to reflect refreshed Market Board prices.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 4458
$ws.Range("I98").Value = 4658.5557
$ws.Range("J98").Value = 848
$ws.Range("K98").Value = 4658.5557
$ws.Range("L98").Value = 848
$ws.Range("M98").Value = -3160.5557
$ws.Range("N98").Value = -3844
$ws.Range("H122").Value = 4458
$ws.Range("I122").Value = 4658.5557
$ws.Range("J122").Value = 848
$ws.Range("K122").Value = 13975.6671
$ws.Range("L122").Value = 2544
$ws.Range("M122").Value = -11525.6671
$ws.Range("N122").Value = -7444
$ws.Range("H124").Value = 65491
$ws.Range("J124").Value = 65491
$ws.Range("L124").Value = 65491
$ws.Range("N124").Value = -75311
$ws.Range("H125").Value = 333333950
$ws.Range("I125").Value = 333333950
$ws.Range("J125").Value = 0
$ws.Range("K125").Value = 3000005550
$ws.Range("L125").Value = 0
$ws.Range("M125").Value = -3000003090
$ws.Range("N125").ClearContents()
$ws.Range("H128").Value = 59333.332
$ws.Range("J128").Value = 62500
$ws.Range("L128").Value = 62500
$ws.Range("N128").Value = -72460
$ws.Range("H132").Value = 1666.7887
$ws.Range("I132").Value = 1646.2059
$ws.Range("K132").Value = 4938.6177
$ws.Range("M132").Value = -2408.6177
$ws.Range("H137").Value = 4243.189
$ws.Range("I137").Value = 4874.9165
$ws.Range("K137").Value = 14624.7495
$ws.Range("M137").Value = -12074.7495
$ws.Range("H138").Value = 4213.673
$ws.Range("J138").Value = 7564.852
$ws.Range("L138").Value = 22694.556
$ws.Range("N138").Value = -32974.556

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 2486.386
$ws.Range("I74").Value = 1793.7727
$ws.Range("K74").Value = 1793.7727
$ws.Range("M74").Value = -919.7727
$ws.Range("H77").Value = 2486.386
$ws.Range("I77").Value = 1793.7727
$ws.Range("K77").Value = 8968.863499999999
$ws.Range("M77").Value = -4600.863499999999
$ws.Range("H102").Value = 1654.2307
$ws.Range("I102").Value = 1654.2307
$ws.Range("K102").Value = 1654.2307
$ws.Range("M102").Value = -32.23070000000007
$ws.Range("H113").Value = 53273
$ws.Range("J113").Value = 53273
$ws.Range("L113").Value = 53273
$ws.Range("N113").Value = -61951
$ws.Range("H131").Value = 59519
$ws.Range("J131").Value = 59519
$ws.Range("L131").Value = 59519
$ws.Range("N131").Value = -69599

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H7").Value = 800
$ws.Range("I7").Value = 800
$ws.Range("K7").Value = 800
$ws.Range("M7").Value = -687
$ws.Range("H86").Value = 79905.46000000001
$ws.Range("I86").Value = 144424.58
$ws.Range("J86").Value = 4633.1665
$ws.Range("K86").Value = 144424.58
$ws.Range("L86").Value = 4633.1665
$ws.Range("M86").Value = -143301.58
$ws.Range("N86").Value = -6879.1665
$ws.Range("H89").Value = 79905.46000000001
$ws.Range("I89").Value = 144424.58
$ws.Range("J89").Value = 4633.1665
$ws.Range("K89").Value = 722122.8999999999
$ws.Range("L89").Value = 23165.8325
$ws.Range("M89").Value = -716506.8999999999
$ws.Range("N89").Value = -34397.8325
$ws.Range("H99").Value = 7577968
$ws.Range("I99").Value = 2254.25
$ws.Range("J99").Value = 22729396
$ws.Range("K99").Value = 2254.25
$ws.Range("L99").Value = 22729396
$ws.Range("M99").Value = -756.25
$ws.Range("N99").Value = -22732392
$ws.Range("H134").Value = 10876940
$ws.Range("J134").Value = 8480.5625
$ws.Range("L134").Value = 25441.6875
$ws.Range("N134").Value = -30511.6875

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 7760.7446
$ws.Range("I31").Value = 3160
$ws.Range("J31").Value = 12169.792
$ws.Range("K31").Value = 3160
$ws.Range("L31").Value = 12169.792
$ws.Range("M31").Value = -2865
$ws.Range("N31").Value = -12759.792
$ws.Range("H34").Value = 7760.7446
$ws.Range("I34").Value = 3160
$ws.Range("J34").Value = 12169.792
$ws.Range("K34").Value = 3160
$ws.Range("L34").Value = 12169.792
$ws.Range("M34").Value = -2958
$ws.Range("N34").Value = -12573.792
$ws.Range("H58").Value = 6761.049
$ws.Range("I58").Value = 2351.5625
$ws.Range("K58").Value = 2351.5625
$ws.Range("M58").Value = -2148.5625
$ws.Range("H86").Value = 15631500
$ws.Range("I86").Value = 15631500
$ws.Range("K86").Value = 15631500
$ws.Range("M86").Value = -15630377
$ws.Range("H89").Value = 15631500
$ws.Range("I89").Value = 15631500
$ws.Range("K89").Value = 78157500
$ws.Range("M89").Value = -78151884
$ws.Range("H99").Value = 6200.1304
$ws.Range("I99").Value = 4659.2
$ws.Range("K99").Value = 4659.2
$ws.Range("M99").Value = -3161.2
$ws.Range("H107").Value = 2200.2856
$ws.Range("I107").Value = 2276.2727
$ws.Range("J107").Value = 2151.1177
$ws.Range("K107").Value = 2276.2727
$ws.Range("L107").Value = 2151.1177
$ws.Range("M107").Value = -356.2727
$ws.Range("N107").Value = -5991.1177
$ws.Range("H126").Value = 6200.1304
$ws.Range("I126").Value = 4659.2
$ws.Range("K126").Value = 13977.6
$ws.Range("M126").Value = -11507.6
$ws.Range("H132").Value = 4907.9287
$ws.Range("I132").Value = 2416.4546
$ws.Range("J132").Value = 7648.55
$ws.Range("K132").Value = 7249.3638
$ws.Range("L132").Value = 22945.65
$ws.Range("M132").Value = -4719.3638
$ws.Range("N132").Value = -28005.65
$ws.Range("H136").Value = 6761.049
$ws.Range("I136").Value = 2351.5625
$ws.Range("K136").Value = 7054.6875
$ws.Range("M136").Value = -4504.6875

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 79266050
$ws.Range("J4").Value = 500050
$ws.Range("L4").Value = 1500150
$ws.Range("N4").Value = -1500374
$ws.Range("H56").Value = 6996.75
$ws.Range("I56").Value = 6996.75
$ws.Range("K56").Value = 6996.75
$ws.Range("M56").Value = -6466.75
$ws.Range("H92").Value = 8548841
$ws.Range("J92").Value = 8548841
$ws.Range("L92").Value = 25646523
$ws.Range("N92").Value = -25649019
$ws.Range("H97").Value = 280.375
$ws.Range("I97").Value = 340.75
$ws.Range("J97").Value = 220
$ws.Range("K97").Value = 1022.25
$ws.Range("L97").Value = 660
$ws.Range("M97").Value = -526.25
$ws.Range("N97").Value = -1652
$ws.Range("H113").Value = 6567.5
$ws.Range("J113").Value = 7078.8335
$ws.Range("L113").Value = 21236.5005
$ws.Range("N113").Value = -25576.5005
$ws.Range("H132").Value = 6128.3335
$ws.Range("J132").Value = 7696.579
$ws.Range("L132").Value = 69269.211
$ws.Range("N132").Value = -74329.211
$ws.Range("H133").Value = 7899.5
$ws.Range("I133").Value = 6333
$ws.Range("J133").Value = 10249.25
$ws.Range("K133").Value = 18999
$ws.Range("L133").Value = 30747.75
$ws.Range("M133").Value = -13939
$ws.Range("N133").Value = -40867.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H57").Value = 59442.062
$ws.Range("J57").Value = 59442.062
$ws.Range("L57").Value = 59442.062
$ws.Range("N57").Value = -61082.062
$ws.Range("H80").Value = 4321.9
$ws.Range("I80").Value = 4259.857
$ws.Range("J80").Value = 4466.6665
$ws.Range("K80").Value = 4259.857
$ws.Range("L80").Value = 4466.6665
$ws.Range("M80").Value = -3261.857
$ws.Range("N80").Value = -6462.6665
$ws.Range("H83").Value = 4321.9
$ws.Range("I83").Value = 4259.857
$ws.Range("J83").Value = 4466.6665
$ws.Range("K83").Value = 21299.285
$ws.Range("L83").Value = 22333.3325
$ws.Range("M83").Value = -16307.285
$ws.Range("N83").Value = -32317.3325
$ws.Range("H124").Value = 75540.5
$ws.Range("J124").Value = 75540.5
$ws.Range("L124").Value = 75540.5
$ws.Range("N124").Value = -85360.5
$ws.Range("H126").Value = 9000
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 9000
$ws.Range("K126").Value = 0
$ws.Range("L126").Value = 27000
$ws.Range("M126").ClearContents()
$ws.Range("N126").Value = -31940

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1936.5883
$ws.Range("I22").Value = 675.7778
$ws.Range("K22").Value = 675.7778
$ws.Range("M22").Value = -380.7778
$ws.Range("H27").Value = 1936.5883
$ws.Range("I27").Value = 675.7778
$ws.Range("K27").Value = 675.7778
$ws.Range("M27").Value = -568.7778
$ws.Range("H55").Value = 494.05
$ws.Range("I55").Value = 186.5
$ws.Range("K55").Value = 186.5
$ws.Range("M55").Value = -13.5
$ws.Range("H68").Value = 8098.4116
$ws.Range("I68").Value = 6576.1113
$ws.Range("J68").Value = 9811
$ws.Range("K68").Value = 6576.1113
$ws.Range("L68").Value = 9811
$ws.Range("M68").Value = -5827.1113
$ws.Range("N68").Value = -11309
$ws.Range("H71").Value = 8098.4116
$ws.Range("I71").Value = 6576.1113
$ws.Range("J71").Value = 9811
$ws.Range("K71").Value = 32880.5565
$ws.Range("L71").Value = 49055
$ws.Range("M71").Value = -29136.5565
$ws.Range("N71").Value = -56543
$ws.Range("H125").Value = 60541
$ws.Range("J125").Value = 60541
$ws.Range("L125").Value = 60541
$ws.Range("N125").Value = -70381
$ws.Range("H127").Value = 59440.25
$ws.Range("J127").Value = 59440.25
$ws.Range("L127").Value = 59440.25
$ws.Range("N127").Value = -69360.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H128").Value = 59519
$ws.Range("J128").Value = 59519
$ws.Range("L128").Value = 59519
$ws.Range("N128").Value = -69479
$ws.Range("H136").Value = 19806130
$ws.Range("I136").Value = 27778600
$ws.Range("K136").Value = 83335800
$ws.Range("M136").Value = -83333250

Write-Host "Applied all Leve profit updates."
